$wb = $excel.ActiveWorkbook

# --- Update selection on "Town of Middletown" (was A8, becomes C2) ---
$wsMiddletown = $wb.Worksheets.Item("Town of Middletown")
$wsMiddletown.Activate() | Out-Null
$wsMiddletown.Range("C2").Select() | Out-Null

# --- Add the new "Town of Poolesville" worksheet as the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Town of Poolesville"

# --- Header row ---
$newSheet.Range("A1").Value = "Zone"
$newSheet.Range("B1").Value = "Zone Abbreviation"
$newSheet.Range("C1").Value = "Issuing Body"
$newSheet.Range("D1").Value = "Zone General Description"

# --- Data rows ---
$newSheet.Range("A2").Value = "Residential Multi-Family"
$newSheet.Range("B2").Value = "PR MUL"
$newSheet.Range("C2").Value = "Poolesville Planning and Zoning Department"
$newSheet.Range("D2").Value = "Residential zoning districts are zones that are developed to house individuals and families. The Residential Multi-Family zone is for residential properties containing multiple family dwelling units like townhouses. "

$newSheet.Range("A3").Value = "Residential 1/3 Acre Lots"
$newSheet.Range("B3").Value = "R 1/3"
$newSheet.Range("C3").Value = "Poolesville Planning and Zoning Department"
$newSheet.Range("D3").Value = "Residential zoning districts are zones that are developed to house individuals and families. The Residential R 1/3 zone is for residential properties at least 1/3 of an acre in size. "

$newSheet.Range("A4").Value = "Residential ½ Acre Lots"
$newSheet.Range("B4").Value = "R ½"
$newSheet.Range("C4").Value = "Poolesville Planning and Zoning Department"
$newSheet.Range("D4").Value = "Residential zoning districts are zones that are developed to house individuals and families. The Residential 1/2 zone is for residential properties at least 1/2 an acre in size. "

$newSheet.Range("A5").Value = "Residential ¾ Acre Lots"
$newSheet.Range("B5").Value = "R ¾"
$newSheet.Range("C5").Value = "Poolesville Planning and Zoning Department"
$newSheet.Range("D5").Value = "Residential zoning districts are zones that are developed to house individuals and families. The Residential 3/4 zone is for residential properties at least 3/4 of an acre in size. "

$newSheet.Range("A6").Value = "Commercial"
$newSheet.Range("B6").Value = "P COMM"
$newSheet.Range("C6").Value = "Poolesville Planning and Zoning Department"
$newSheet.Range("D6").Value = "Commercial zones are for businesses, restaurants, retail stores, convenience stores, entertainment stores, or automobile shops; however, they also support mixed uses or residential projects. The purpose of the commercial zone is to create a vibrant Town center that serves as a destination for residents and visitors to walk, shop, dine, live, and interact. Development and redevelopment projects should be harmonious with and enhance the characteristics of the Town Center."

# --- Column widths to match source formatting (bestFit-style autosize) ---
$newSheet.Columns.Item(1).ColumnWidth = 21.877604166666668
$newSheet.Columns.Item(2).ColumnWidth = 16.166666666666668
$newSheet.Columns.Item(3).ColumnWidth = 40.022135416666664
$newSheet.Columns.Item(4).ColumnWidth = 254.87760416666666

# --- Make the new sheet the active tab/selection ---
$newSheet.Activate() | Out-Null
$newSheet.Range("B12").Select() | Out-Null
